$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply the numeric format (adds a new cellXfs entry, fillId=4/borderId=1 kept,
# numFmtId 164 -> 165) to the whole "load-times" block (data rows + total row).
$ws.Range("B24:F30").NumberFormat = "0.0000"

# Row 24 - CS1
$ws.Range("B24").Formula = "=(0.148542+0.148461+0.147445+0.148781+0.148514)/5"
$ws.Range("C24").Formula = "=(0.22119+0.195683+0.196896+0.195851+0.200137)/5"
$ws.Range("D24").Formula = "=(1.02987+1.04242+1.02666+1.0302+1.03145)/5"
$ws.Range("E24").Formula = "=(15.0838+15.8622+15.2584+15.4685+15.2169)/5"
$ws.Range("F24").Formula = "=(1.07828+1.05154+1.03035+1.03241+1.02617)/5"

# Row 25 - PA1
$ws.Range("B25").Formula = "=(0.656544+0.656577+0.656983+0.657071+0.657465)/5"
$ws.Range("C25").Formula = "=(0.762904+0.762778+0.763616+0.763524+0.764375)/5"
$ws.Range("D25").Formula = "=(3.61774+3.61151+3.61306+3.65721+3.64153)/5"
$ws.Range("E25").Formula = "=(218.66+218.778+221.719+221.025+218.572)/5"
$ws.Range("F25").Formula = "=(3.62415+3.61532+3.61154+3.6341+3.61377)/5"

# Row 26 - SP1
$ws.Range("B26").Formula = "=(3.07683+3.06376+3.13867+3.07326+3.06826)/5"
$ws.Range("C26").Formula = "=(3.84828+3.85047+3.85387+3.85881+3.85288)/5"
$ws.Range("D26").Formula = "=(22.8852+22.8723+23.049+24.2356+22.8745)/5"
$ws.Range("E26").Formula = "=(700.848+694.558+690.941+714.705+693.824)/5"
$ws.Range("F26").Formula = "=(23.0711+22.9609+22.9109+24.6303+22.8736)/5"

# Row 27 - SP2
$ws.Range("B27").Formula = "=(11.3389+11.4321+11.58+11.3597+11.3709)/5"
$ws.Range("C27").Formula = "=(13.3618+13.2946+13.3055+13.3022+13.301)/5"
$ws.Range("D27").Formula = "=(53.0528+53.2677+53.8663+53.3803+53.6128)/5"
$ws.Range("E27").Formula = "=(3231.63+3276.54+3198.06+3075.7+2998.18)/5"
$ws.Range("F27").Formula = "=(662.312+397.023+349.231+398.163+349.231+371.301)/5"

# Row 28 - PA2
$ws.Range("B28").Formula = "=(2.25127+2.25606+2.4884+2.25569+2.25729)/5"
$ws.Range("C28").Formula = "=(2.90764+2.90838+2.90881+2.90827+2.91081)/5"
$ws.Range("D28").Formula = "=(18.2144+18.2161+18.7175+18.3709+18.2204)/5"
$ws.Range("E28").Formula = "=(36.996+39.076+36.1302+34.3196+37.0837)/5"
$ws.Range("F28").Formula = "=(20.5763+20.6189+20.6514+20.8378+20.6164)/5"

# Row 29 - CS2
$ws.Range("B29").Formula = "=(0.16176+0.161215+0.162113+0.16176+0.16247)/5"
$ws.Range("C29").Formula = "=(0.214602+0.214146+0.214403+0.214654+0.214532)/5"
$ws.Range("D29").Formula = "=(1.145+1.15569+1.14106+1.14331+1.14345)/5"
$ws.Range("E29").Formula = "=(13.1867+13.345+13.3665+13.1004+13.1382)/5"
$ws.Range("F29").Formula = "=(1.14281+1.14335+1.14226+1.14081+1.14028)/5"

# Row 30 - Total
$ws.Range("B30").Formula = "=SUM(B24:B29)"
$ws.Range("C30").Formula = "=SUM(C24:C29)"
$ws.Range("D30").Formula = "=SUM(D24:D29)"
$ws.Range("E30").Formula = "=SUM(E24:E29)"
$ws.Range("F30").Formula = "=SUM(F24:F29)"

# Move the active selection to match the author's final cursor position.
[void]$ws.Range("F28").Select()
